$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.525965269531893
$ws.Range("C2").Value = 1.322619210585525
$ws.Range("D2").Value = 1.169002513005085
$ws.Range("E2").Value = 1.220884036319589
$ws.Range("F2").Value = 1.332915022613899
$ws.Range("G2").Value = 1.314277210411198

$ws.Range("B3").Value = 1.539836047058403
$ws.Range("C3").Value = 1.24644318523582
$ws.Range("D3").Value = 1.236655072989541
$ws.Range("E3").Value = 1.202193997384603
$ws.Range("F3").Value = 1.189771017454148
$ws.Range("G3").Value = 1.282979864024503

$ws.Range("B4").Value = 1.41628701880448
$ws.Range("C4").Value = 1.060330398414768
$ws.Range("D4").Value = 1.269980323660571
$ws.Range("E4").Value = 1.239977132725892
$ws.Range("F4").Value = 1.294330520127511
$ws.Range("G4").Value = 1.256181078746645

$ws.Range("B5").Value = 1.54346131258994
$ws.Range("C5").Value = 1.175073221374955
$ws.Range("D5").Value = 0.978708444294902
$ws.Range("E5").Value = 0.9916108842675304
$ws.Range("F5").Value = 1.14343627414154
$ws.Range("G5").Value = 1.166458027333773

$ws.Range("B6").Value = 0.9240847822018921
$ws.Range("C6").Value = 0.7970027915421876
$ws.Range("D6").Value = 0.5682791430575371
$ws.Range("E6").Value = 0.5827491581979505
$ws.Range("F6").Value = 0.707176145132798
$ws.Range("G6").Value = 0.7158584040264729

$ws.Range("B7").Value = 0.3957573216720039
$ws.Range("C7").Value = 0.5660907733492631
$ws.Range("D7").Value = 0.3794973895153716
$ws.Range("E7").Value = 0.3192274030934538
$ws.Range("F7").Value = 0.4301859734429974
$ws.Range("G7").Value = 0.4181517722146179

$ws.Range("B8").Value = 0.3760537054939977
$ws.Range("C8").Value = 0.3420546008335217
$ws.Range("D8").Value = 0.1578135723174616
$ws.Range("E8").Value = 0.1412350080299828
$ws.Range("F8").Value = 0.2201480844511626
$ws.Range("G8").Value = 0.2474609942252253

$ws.Range("B9").Value = 0.7108270233441099
$ws.Range("C9").Value = 0.2090492859052139
$ws.Range("D9").Value = 0.09728488920783651
$ws.Range("E9").Value = 0.09087548139238655
$ws.Range("F9").Value = 0.1282813427818635
$ws.Range("G9").Value = 0.2472636045262821

$ws.Range("B10").Value = 0.9788722036769728
$ws.Range("C10").Value = 0.2096266312732893
$ws.Range("D10").Value = 0.09393413116731208
$ws.Range("E10").Value = 0.05070654976853972
$ws.Range("F10").Value = 0.1487939289094512
$ws.Range("G10").Value = 0.296386688959113

$ws.Range("B11").Value = 1.899489556631275
$ws.Range("C11").Value = 0.1982595958728413
$ws.Range("D11").Value = 0.1169534596228009
$ws.Range("E11").Value = 0.04180567673729096
$ws.Range("F11").Value = 0.1259859422456563
$ws.Range("G11").Value = 0.4764988462219729

$ws.Range("B12").Value = 1.536941295366478
$ws.Range("C12").Value = 0.1932361358617809
$ws.Range("D12").Value = 0.1370453144309562
$ws.Range("E12").Value = 0.04742613787190875
$ws.Range("F12").Value = 0.1408756989364836
$ws.Range("G12").Value = 0.4111049164935216

$ws.Range("B13").Value = 2.149965669878457
$ws.Range("C13").Value = 0.2020063127782311
$ws.Range("D13").Value = 0.1422063693269555
$ws.Range("E13").Value = 0.05611950484942585
$ws.Range("F13").Value = 0.1479585045290324
$ws.Range("G13").Value = 0.5396512722724205

$ws.Range("B14").Value = 1.558743048855251
$ws.Range("C14").Value = 0.2059833804856536
$ws.Range("D14").Value = 0.1547725216578992
$ws.Range("E14").Value = 0.06183024005044346
$ws.Range("F14").Value = 0.1481775129269887
$ws.Range("G14").Value = 0.4259013407952471

$ws.Range("B15").Value = 1.853629939407736
$ws.Range("C15").Value = 0.2053231409355986
$ws.Range("D15").Value = 0.1677313618360533
$ws.Range("E15").Value = 0.0599667715278106
$ws.Range("F15").Value = 0.1477482588111194
$ws.Range("G15").Value = 0.4868798945036635

$ws.Range("B16").Value = 1.789704015869725
$ws.Range("C16").Value = 0.2048315381960227
$ws.Range("D16").Value = 0.1694420174567767
$ws.Range("E16").Value = 0.06445625548496803
$ws.Range("F16").Value = 0.1502687465957444
$ws.Range("G16").Value = 0.4757405147206473

$ws.Range("B17").Value = 1.801132914686453
$ws.Range("C17").Value = 0.2050245889580702
$ws.Range("D17").Value = 0.171955657599162
$ws.Range("E17").Value = 0.06326199470326921
$ws.Range("F17").Value = 0.151041449405909
$ws.Range("G17").Value = 0.4784833210705727

$ws.Range("B18").Value = 1.799097682926146
$ws.Range("C18").Value = 0.2053059878068844
$ws.Range("D18").Value = 0.1720745515957336
$ws.Range("E18").Value = 0.06356284793215772
$ws.Range("F18").Value = 0.1508856369303
$ws.Range("G18").Value = 0.4781853414382443

